$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.03375954870909267
$ws.Range("C2").Value = 0.6485414634705653
$ws.Range("D2").Value = 0.7314436294522902
$ws.Range("E2").Value = 0.8552447775065862
$ws.Range("F2").Value = 0.8663660556530268
$ws.Range("G2").Value = 37

$ws.Range("B3").Value = 0.03135085916488722
$ws.Range("C3").Value = 0.5511937747325102
$ws.Range("D3").Value = 0.5928526553971856
$ws.Range("E3").Value = 0.7699692561376626
$ws.Range("F3").Value = 0.7802437722544544
$ws.Range("G3").Value = 36

$ws.Range("B4").Value = 0.0324755997945022
$ws.Range("C4").Value = 0.5741707009966609
$ws.Range("D4").Value = 0.6201381699606496
$ws.Range("E4").Value = 0.7874885205262675
$ws.Range("F4").Value = 0.7983056080049333
$ws.Range("G4").Value = 35

$ws.Range("B5").Value = 0.04451789004105118
$ws.Range("C5").Value = 0.5528601565202774
$ws.Range("D5").Value = 0.615403550752249
$ws.Range("E5").Value = 0.7844766094360296
$ws.Range("F5").Value = 0.7949907199654754
$ws.Range("G5").Value = 34

$ws.Range("B6").Value = 0.02488288622478917
$ws.Range("C6").Value = 0.5597336268738891
$ws.Range("D6").Value = 0.6071832843101319
$ws.Range("E6").Value = 0.7792196637085924
$ws.Range("F6").Value = 0.7908977527023375
$ws.Range("G6").Value = 33

$ws.Range("B7").Value = 0.06447809342936392
$ws.Range("C7").Value = 0.5639312790710509
$ws.Range("D7").Value = 0.6386050102870782
$ws.Range("E7").Value = 0.7991276558141873
$ws.Range("F7").Value = 0.8092673457567485
$ws.Range("G7").Value = 32

$ws.Range("B8").Value = 0.0316477260114686
$ws.Range("C8").Value = 0.5875569779702665
$ws.Range("D8").Value = 0.6486140808091702
$ws.Range("E8").Value = 0.8053658055872314
$ws.Range("F8").Value = 0.8180462001963719
$ws.Range("G8").Value = 31

$ws.Range("B9").Value = 0.05575202370392403
$ws.Range("C9").Value = 0.5712956843590213
$ws.Range("D9").Value = 0.6777300950868858
$ws.Range("E9").Value = 0.8232436426033825
$ws.Range("F9").Value = 0.8353948933700528
$ws.Range("G9").Value = 30

$ws.Range("B10").Value = 0.03351674390117437
$ws.Range("C10").Value = 0.5951644878936757
$ws.Range("D10").Value = 0.6821752870341646
$ws.Range("E10").Value = 0.8259390334850174
$ws.Range("F10").Value = 0.8398661783807415
$ws.Range("G10").Value = 29

$ws.Range("B11").Value = 0.06637440518007821
$ws.Range("C11").Value = 0.5683437572335094
$ws.Range("D11").Value = 0.6866267435578414
$ws.Range("E11").Value = 0.8286294368159035
$ws.Range("F11").Value = 0.8411234351010097
$ws.Range("G11").Value = 28

